$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 26, shifting existing rows 26-62 down to 27-63.
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with the new weekly price entry.
$ws.Cells.Item(26, 1).Value = 11
$ws.Cells.Item(26, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(26, 3).Value = "Bíobío"
$ws.Cells.Item(26, 4).Value = 45240
$ws.Cells.Item(26, 5).Value = 8
$ws.Cells.Item(26, 6).Value = "Fruta"
$ws.Cells.Item(26, 7).Value = 100107
$ws.Cells.Item(26, 8).Value = "Otros"
$ws.Cells.Item(26, 9).Value = 100107002
$ws.Cells.Item(26, 10).Value = "Chirimoya"
$ws.Cells.Item(26, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(26, 12).Value = "Primera"
$ws.Cells.Item(26, 13).Value = 250
$ws.Cells.Item(26, 14).Value = 22000
$ws.Cells.Item(26, 15).Value = 23000
$ws.Cells.Item(26, 16).Value = 22400
$ws.Cells.Item(26, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(26, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(26, 19).Value = 2240
$ws.Cells.Item(26, 20).Value = 10
